# Corrected some selection scopes
#
# The matched-errors table had its first data row (the 2020-04-01 / "Q1"
# quarter) missing. This inserts a new row at row 3 (pushing the existing
# rows 3-22 down to 4-23, which is what the target diff shows: every row's
# B:H content now equals what used to be one row above it, and a brand new
# row of data appears right under the "2020-01-01" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 3:22 down to 4:23, inserting a blank row 3.
$ws.Rows("3:3").Insert()

# The inserted row has no formatting yet; clone it from the row above (A2),
# which carries the bordered/bold/centered label style used by every entry
# in column A, so the new label cell matches its neighbours.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new quarter's label and its 7 (previously-missing) data points.
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"
$ws.Range("B3").Value = 9.643547872076862
$ws.Range("C3").Value = -8.527713928060606
$ws.Range("D3").Value = -0.388042381081458
$ws.Range("E3").Value = 1.759528090717934
$ws.Range("F3").Value = -1.680501895805395
$ws.Range("G3").Value = -1.741676074219596
$ws.Range("H3").Value = 0.2997798629366579
